{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Annual Report\" paragraph and the \"Event Report\" paragraph.\nlet annualPara = null;\nlet eventPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Annual Report\") !== -1) {\n    annualPara = p;\n  } else if (p.text.indexOf(\"Event Report\") !== -1) {\n    eventPara = p;\n  }\n}\n\n// Remove the existing \"_GoBack\" bookmark (it currently sits at the end of\n// the \"Event Report\" paragraph) before the paragraph layout shifts.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Delete the whole \"Annual Report\" assumption paragraph.\nif (annualPara) {\n  annualPara.delete();\n}\nawait context.sync();\n\n// Re-insert the \"_GoBack\" bookmark at the very start of the \"Event Report\"\n// paragraph, matching where Word leaves the last-edit marker after the\n// preceding paragraph was removed.\nif (eventPara) {\n  eventPara.getRange(\"Start\").insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last-edit marker (\"_GoBack\") currently sits at the end of the \"Event\n# Report\" paragraph. Drop it now -- before paragraph indices shift -- we'll\n# re-add it in its new spot once the \"Annual Report\" paragraph is gone.\ntry {\n    $d.Bookmarks(\"_GoBack\").Delete()\n} catch {\n}\n\n# Locate and remove the whole \"Annual Report\" assumption paragraph.\n$findRange = $d.Content\n$null = $findRange.Find.Execute(\"Annual Report\")\n$annualPara = $findRange.Paragraphs(1)\n$annualPara.Range.Delete()\n\n# Re-insert \"_GoBack\" at the very start of the \"Event Report\" paragraph --\n# where Word leaves the last-edit marker once the preceding paragraph is\n# removed.\n$findRange2 = $d.Content\n$null = $findRange2.Find.Execute(\"Event Report\")\n$eventPara = $findRange2.Paragraphs(1)\n$target = $eventPara.Range.Duplicate\n$target.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $target)\n"}
